$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the "Objetivos:" (row 10) value — it was incorrectly holding the
#    professor's name; replace with the real Portuguese objectives text.
$ws.Range("B10").Value = "Apresentar ao estudante conceitos gerais sobre a síntese de polímeros, destacando as principais vias usadas para a geração de materiais poliméricos na atualidade. Capacitar o estudante para relacionar a síntese com a estrutura, o comportamento e a utilização de polímeros."
$ws.Range("C10").Value = "Apresentar ao estudante conceitos gerais sobre a síntese de polímeros, destacando as principais vias usadas para a geração de materiais poliméricos na atualidade. Capacitar o estudante para relacionar a síntese com a estrutura, o comportamento e a utilização de polímeros."

# 2. Insert a new row at 13 for "Docentes responsaveis:" value
#    (shifts everything from old row 13 down to row 14, etc.)
$ws.Rows("13:13").Insert()
$ws.Range("A13").Clear()

# 3. Populate the new row 13 (B/C) with the professor's name, copying
#    formatting (style) from a correctly-styled neighboring value cell.
$ws.Range("B19").Copy($ws.Range("B13"))
$ws.Range("C19").Copy($ws.Range("C13"))
$ws.Range("B13").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C13").Value = "1033242 - Fábio Herbst Florenzano"

# 4. "Programa resumido:" value (now row 14) — replace placeholder "Semestral"
$ws.Range("B14").Value = "Fundamentos: massa molar média e conceitos gerais sobre química de polímeros. Polímeros de condensação e de adição. Polimerização em cadeia. Polimerização em etapas. Métodos de Polimerização. Modificação e degradação de polímeros e seu retardo."
$ws.Range("C14").Value = "Fundamentos: massa molar média e conceitos gerais sobre química de polímeros. Polímeros de condensação e de adição. Polimerização em cadeia. Polimerização em etapas. Métodos de Polimerização. Modificação e degradação de polímeros e seu retardo."

# 5. "Programa:" value (now row 16) — replace misplaced date with real syllabus
$ws.Range("B16").Value = "Fundamentos: massa molar média e conceitos gerais sobre química de polímeros. Polímeros de condensação e de adição: conceitos gerais. Polimerização em Etapas. Polimerização em cadeia: via radical, aniônica e catiônica. Polimerização via desativação reversível de radicais. Polimerização estéreo-específica: Ziegler-Natta e outras. Métodos de polimerização: batelada, solução, suspensão, emulsão e interfacial. Modificação de polímeros: reticulação e vulcanização; modificações em aromáticos, hidrólise e outras. Química da degradação de polímeros: processos gerais e métodos de controle."
$ws.Range("C16").Value = "Fundamentos: massa molar média e conceitos gerais sobre química de polímeros. Polímeros de condensação e de adição: conceitos gerais. Polimerização em Etapas. Polimerização em cadeia: via radical, aniônica e catiônica. Polimerização via desativação reversível de radicais. Polimerização estéreo-específica: Ziegler-Natta e outras. Métodos de polimerização: batelada, solução, suspensão, emulsão e interfacial. Modificação de polímeros: reticulação e vulcanização; modificações em aromáticos, hidrólise e outras. Química da degradação de polímeros: processos gerais e métodos de controle."

# 6. "Metodo:" value (now row 19) — was holding the professor's name, fix it
$ws.Range("B19").Value = "Provas escritas envolvendo o conteúdo teórico ministrado em sala de aula."
$ws.Range("C19").Value = "Provas escritas envolvendo o conteúdo teórico ministrado em sala de aula."

# 7. "Criterio:" value (now row 20) — shift text that used to be here down
$ws.Range("B20").Value = "Duas avaliações, sendo que a nota final corresponde à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados. Alunos com notas finais situadas no intervalo de 3 a 5 serão encaminhados à recuperação."
$ws.Range("C20").Value = "Duas avaliações, sendo que a nota final corresponde à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados. Alunos com notas finais situadas no intervalo de 3 a 5 serão encaminhados à recuperação."

# 8. "Norma de recuperacao:" value (now row 21)
$ws.Range("B21").Value = "O aluno será submetido a um programa de estudos destinado a rever o conteúdo abordado na disciplina. Ao final deste período será aplicada uma nova avaliação. A nota final do aluno será a média aritmética desta avaliação com a nota anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5."
$ws.Range("C21").Value = "O aluno será submetido a um programa de estudos destinado a rever o conteúdo abordado na disciplina. Ao final deste período será aplicada uma nova avaliação. A nota final do aluno será a média aritmética desta avaliação com a nota anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5."

# 9. "Bibliografia:" value (now row 22) — replace with the real bibliography
$ws.Range("B22").Value = "G. ODIAN Principles of Polymerization, 3rd Edition, New York: Wiley-Interscience, 1991.`nF. W. Billmeyer. Textbook of Polymer Chemistry, 3rd edition, New York: Wiley-Interscience, 1984.`nC. E. Carraher. Introduction to Polymer Chemistry, 1st Edition, Boca Raton: Taylor and Francis, 2010.`nS. V. Canevarolo. Ciência dos Polímeros: um texto básico para Engenheiros e Tecnólogos, 2ª. edição, São Paulo: Artliber, 2006."
$ws.Range("C22").Value = "G. ODIAN Principles of Polymerization, 3rd Edition, New York: Wiley-Interscience, 1991.`nF. W. Billmeyer. Textbook of Polymer Chemistry, 3rd edition, New York: Wiley-Interscience, 1984.`nC. E. Carraher. Introduction to Polymer Chemistry, 1st Edition, Boca Raton: Taylor and Francis, 2010.`nS. V. Canevarolo. Ciência dos Polímeros: um texto básico para Engenheiros e Tecnólogos, 2ª. edição, São Paulo: Artliber, 2006."
